$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# totalTime changed from 1200 to 15
$ws.Range("C5").Value = 15

# Insert a new row at 10 for "stopWhenDead" (shifts old rows 11+ down by one)
$ws.Rows("10:10").Insert()

# Insert a new row at 54 for "collisionKillDistance" (shifts old rows 54+ down by two total)
$ws.Rows("54:54").Insert()

# Fill in the new collisionKillDistance row first so its shared string is
# registered before stopWhenDead's (matches authoring order)
$ws.Range("A54").Value = "[m]"
$ws.Range("B54").Value = "collisionKillDistance"
$ws.Range("C54").Value = 4

$ws.Range("B10").Value = "stopWhenDead"
$ws.Range("C10").Value = $true

# Update the selected/active cell shown when the sheet is reopened
[void]$ws.Range("C8").Select()
